# Fixing the big mistake: correct Total (B) and Community (D) monthly
# consumption values on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = 3900.570070066671;  D = 263.3084614166667 }
    3  = @{ B = 3684.504820050005;  D = 240.0477465333333 }
    4  = @{ B = 3935.517170283338;  D = 277.7661802 }
    5  = @{ B = 3792.569330550004;  D = 256.7541351166667 }
    6  = @{ B = 3941.616397850005;  D = 271.8358439 }
    7  = @{ B = 3805.007706733338;  D = 260.4824763666667 }
    8  = @{ B = 3939.895343700005;  D = 267.8358084333333 }
    9  = @{ B = 3925.973884166671;  D = 263.1251272666667 }
    10 = @{ B = 3794.700280216671;  D = 257.78580565 }
    11 = @{ B = 3926.574312566671;  D = 264.3467976833333 }
    12 = @{ B = 3789.640982833338;  D = 261.13414755 }
    13 = @{ B = 3770.267920050005;  D = 258.2821367333333 }
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row].B
    $ws.Range("D$row").Value = $values[$row].D
}
